$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was switched to Page Layout view and the user clicked into cell
# F15 (the first empty row below the table) and centered it vertically,
# which is enough to create a new (otherwise blank) row 15 with a
# vertical-center style.
$ws.Range("F15").VerticalAlignment = -4108

# Page Layout view surfaces the header/footer editing regions; the user
# filled in their name / surname / student number.
$ws.PageSetup.LeftHeader = "Name: Thozamile "
$ws.PageSetup.CenterHeader = "Surname: Madela"
$ws.PageSetup.RightHeader = "ITS Number: 202411681`n"

# The print scale was adjusted so the report fits the page better.
$ws.PageSetup.Zoom = 63

# Page Layout view re-flows drawing anchors slightly against the page
# margins/header area, nudging the chart down and to the right a touch
# while keeping its size identical.
$co = $ws.ChartObjects(1)
$co.Left = $co.Left + 4.5
$co.Top = $co.Top + 3.75

# Leave the selection on the newly touched cell, matching the saved
# workbook's last cursor position.
$null = $ws.Range("F15").Select()
